$d = $word.ActiveDocument

function Set-RunTextWithStyle($range, $text, $rStyle) {
    # Replace a range's contents with a single run (optionally styled) whose
    # text carries xml:space="preserve", matching this document's authoring
    # convention. Uses a Flat-OPC wrapped InsertXML so run formatting
    # (w:rPr/w:rStyle) is fully controlled and neighboring runs are left
    # untouched (avoids accidental run-coalescing from plain text edits).
    $rpr = ''
    if ($rStyle) {
        $rpr = '<w:rPr><w:rStyle w:val="' + $rStyle + '"/></w:rPr>'
    }
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body><w:p><w:r>' + $rpr + '<w:t xml:space="preserve">' + $text + '</w:t></w:r></w:p></w:body>' + `
        '</w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($xml)
}

# 1) "Gopi Santhiran ()" -> "Gopi Santhiran (1001371534)"
#    Only the "()" run's content changes; replace that run's range in place
#    so the preceding " " run is left intact (not coalesced into it).
$pAuthor = $d.Paragraphs(5)
$rAuthor = $pAuthor.Range
$authorStart = $rAuthor.Start
$authorText = $rAuthor.Text
$parenOffset = $authorText.IndexOf("()")
$parenRange = $d.Range($authorStart + $parenOffset, $authorStart + $parenOffset + 2)
Set-RunTextWithStyle $parenRange "(1001371534)" $null

# 2) Remove the whole paragraph that holds the inline figure
#    (unnamed-chunk-4-1.png), style "FirstParagraph".
$pFigure = $d.Paragraphs(51)
$pFigure.Range.Delete()

# 3) "## [1] 162 182" -> "## [1] 803.619"
$pOutliers = $d.Paragraphs(51)
$rOutliers = $pOutliers.Range
$outliersSub = $d.Range($rOutliers.Start, $rOutliers.End - 1)
Set-RunTextWithStyle $outliersSub "## [1] 803.619" "VerbatimChar"

# 4) Collapse the two-run/line-break paragraph into a single run.
$pRstudent = $d.Paragraphs(52)
$rRstudent = $pRstudent.Range
$rstudentSub = $d.Range($rRstudent.Start, $rRstudent.End - 1)
Set-RunTextWithStyle $rstudentSub "## [1] 832.5803" "VerbatimChar"
